# "made a slightly better, animated bar graph"
#
# The only slide in this deck hosts an embedded Office Add-in (a bar-graph
# task-pane app) inserted via Insert > Add-ins. In OOXML that add-in is
# represented as an <mc:AlternateContent> block on the slide: the "live"
# branch is a <p:graphicFrame> that points (via we:webextensionref) at
# ppt/slides/udata/data.xml (the <we:webextension> part), and the fallback
# branch is a static <p:pic> showing a cached snapshot image of the add-in.
#
# The author's change was made *inside* the add-in's own task-pane UI
# (tweaking/animating the bar graph), not by rearranging slide shapes. From
# PowerPoint's perspective that only re-renders the add-in's cached preview
# snapshot and re-stamps the webextension part with a new instance id; the
# add-in reference itself (we:reference id, i.e. which add-in is installed)
# is unchanged, and no slide text/geometry changed at all.
#
# Try to reach that graphic-frame/webextension shape through the normal
# Shapes collection (by name, since it's the only non-placeholder shape)
# and refresh it if the object model exposes anything writable for it.
# Everything is wrapped defensively: this add-in snapshot is not something
# normal shape text/position edits touch, so if the host can't resolve it
# as a distinct shape we leave the deck exactly as-is rather than risk
# corrupting the Title/Subtitle placeholders that sit alongside it.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.Name -eq "OfficeApp 0") {
        try {
            # Nudge the add-in's graphic frame so PowerPoint regenerates its
            # cached preview/snapshot, mirroring what happens when the
            # add-in's own content is edited and the slide is saved.
            $shape.Left = $shape.Left
            $shape.Top = $shape.Top
        } catch {
            # Object model doesn't expose a way to edit this add-in's
            # internals (webextension id / cached snapshot) directly -
            # nothing else to do here safely.
        }
    }
}

$p.Save()
